$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0183284457478006
$ws.Range("C2").Value = 0.0227272727272727
$ws.Range("D2").Value = 0.934017595307918
$ws.Range("E2").Value = 0.0234604105571848
$ws.Range("F2").Value = 0.977272727272727
$ws.Range("G2").Value = 0.0205278592375367
$ws.Range("H2").Value = 0.000733137829912023
$ws.Range("I2").Value = 0.0381231671554252
$ws.Range("J2").Value = 0.00586510263929619
$ws.Range("K2").Value = 0.74266862170088
$ws.Range("L2").Value = 0.00219941348973607
$ws.Range("M2").Value = 0.254398826979472
$ws.Range("N2").Value = 0.976539589442815
$ws.Range("O2").Value = 0.0410557184750733
$ws.Range("P2").Value = 0.000733137829912023
$ws.Range("R2").Value = 0.3841642228739
$ws.Range("S2").Value = 0.0183284457478006
$ws.Range("T2").Value = 0.0175953079178886
$ws.Range("U2").Value = 0.960410557184751
$ws.Range("V2").Value = 0.0285923753665689
$ws.Range("W2").Value = 0.030058651026393
$ws.Range("X2").Value = 0.00659824046920821
$ws.Range("B3").Value = 0.113636363636364
$ws.Range("C3").Value = 0.00366568914956012
$ws.Range("D3").Value = 0.000733137829912023
$ws.Range("E3").Value = 0.947947214076246
$ws.Range("F3").Value = 0.0168621700879765
$ws.Range("G3").Value = 0.928152492668622
$ws.Range("H3").Value = 0.971407624633431
$ws.Range("I3").Value = 0.959677419354839
$ws.Range("J3").Value = 0.063049853372434
$ws.Range("K3").Value = 0.0249266862170088
$ws.Range("L3").Value = 0.0161290322580645
$ws.Range("M3").Value = 0.00146627565982405
$ws.Range("N3").Value = 0.00513196480938416
$ws.Range("O3").Value = 0.00219941348973607
$ws.Range("P3").Value = 0.991935483870968
$ws.Range("Q3").Value = 0.997800586510264
$ws.Range("R3").Value = 0.000733137829912023
$ws.Range("S3").Value = 0.00293255131964809
$ws.Range("T3").Value = 0.961143695014663
$ws.Range("U3").Value = 0.00146627565982405
$ws.Range("V3").Value = 0.0102639296187683
$ws.Range("W3").Value = 0.00219941348973607
$ws.Range("X3").Value = 0.0168621700879765
$ws.Range("B4").Value = 0.00366568914956012
$ws.Range("C4").Value = 0.0168621700879765
$ws.Range("D4").Value = 0.06158357771261
$ws.Range("E4").Value = 0.00659824046920821
$ws.Range("F4").Value = 0.00293255131964809
$ws.Range("G4").Value = 0.000733137829912023
$ws.Range("H4").Value = 0.0219941348973607
$ws.Range("I4").Value = 0.00146627565982405
$ws.Range("J4").Value = 0.0256598240469208
$ws.Range("K4").Value = 0.202346041055718
$ws.Range("L4").Value = 0.00293255131964809
$ws.Range("M4").Value = 0.73900293255132
$ws.Range("N4").Value = 0.00586510263929619
$ws.Range("O4").Value = 0.956744868035191
$ws.Range("Q4").Value = 0.000733137829912023
$ws.Range("R4").Value = 0.00219941348973607
$ws.Range("S4").Value = 0.977272727272727
$ws.Range("T4").Value = 0.0175953079178886
$ws.Range("U4").Value = 0.032258064516129
$ws.Range("V4").Value = 0.0124633431085044
$ws.Range("W4").Value = 0.963343108504399
$ws.Range("X4").Value = 0.971407624633431
$ws.Range("B5").Value = 0.863636363636364
$ws.Range("C5").Value = 0.951612903225806
$ws.Range("D5").Value = 0.00366568914956012
$ws.Range("E5").Value = 0.0219941348973607
$ws.Range("F5").Value = 0.00293255131964809
$ws.Range("G5").Value = 0.0498533724340176
$ws.Range("H5").Value = 0.00586510263929619
$ws.Range("I5").Value = 0.000733137829912023
$ws.Range("J5").Value = 0.905425219941349
$ws.Range("K5").Value = 0.030058651026393
$ws.Range("L5").Value = 0.978739002932551
$ws.Range("M5").Value = 0.00366568914956012
$ws.Range("N5").Value = 0.0124633431085044
$ws.Range("P5").Value = 0.00733137829912024
$ws.Range("Q5").Value = 0.00146627565982405
$ws.Range("R5").Value = 0.612903225806452
$ws.Range("S5").Value = 0.000733137829912023
$ws.Range("T5").Value = 0.00366568914956012
$ws.Range("U5").Value = 0.00586510263929619
$ws.Range("V5").Value = 0.948680351906158
$ws.Range("W5").Value = 0.00439882697947214
$ws.Range("X5").Value = 0.00513196480938416
